# Fruta / hortaliza, semanal
# Insert a new weekly price row at row 122 (pushing existing rows 122-188
# down to 123-189) and populate it with the new week's data, mirroring the
# variety/quality/origin of the (now-shifted) row directly below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 122..188 down to 123..189, duplicating formatting of row 122
# into the freshly inserted row.
$ws.Rows(122).Insert()

# Populate the newly inserted row 122 with this week's record.
$ws.Cells.Item(122, 1).Value = 1
$ws.Cells.Item(122, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(122, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(122, 4).Value = 44596
$ws.Cells.Item(122, 5).Value = 15
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100108
$ws.Cells.Item(122, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(122, 9).Value = 100108006
$ws.Cells.Item(122, 10).Value = "Plátano"
$ws.Cells.Item(122, 11).Value = "Sin especificar"
$ws.Cells.Item(122, 12).Value = "Pintón"
$ws.Cells.Item(122, 13).Value = 120
$ws.Cells.Item(122, 14).Value = 17000
$ws.Cells.Item(122, 15).Value = 18000
$ws.Cells.Item(122, 16).Value = 17500
$ws.Cells.Item(122, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(122, 18).Value = "Ecuador"
$ws.Cells.Item(122, 19).Value = 875
$ws.Cells.Item(122, 20).Value = 20
